$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.406.09'
$ws.Range("E2").Value = '  +0.70%  '
$ws.Range("D3").Value = '1.637.71'
$ws.Range("E3").Value = '  +2.16%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D6").Value = '''304.52'
$ws.Range("E6").Value = '  +0.34%  '
$ws.Range("D7").Value = '''0.3727'
$ws.Range("E7").Value = '  -1.27%  '
$ws.Range("D8").Value = '''52.02'
$ws.Range("E8").Value = '  +0.34%  '
$ws.Range("D9").Value = '''0.3615'
$ws.Range("E9").Value = '  -0.74%  '
$ws.Range("D10").Value = '''1.243'
$ws.Range("E10").Value = '  -2.54%  '
$ws.Range("D11").Value = '''0.08096'
$ws.Range("E11").Value = '  -0.61%  '
$ws.Range("E12").Value = '  +0.10%  '
$ws.Range("D13").Value = '''22.76'
$ws.Range("E13").Value = '  -0.62%  '
$ws.Range("D14").Value = '''6.578'
$ws.Range("E14").Value = '  -0.49%  '
$ws.Range("D15").Value = '''0.00001265'
$ws.Range("E15").Value = '  +1.09%  '
$ws.Range("D16").Value = '''7.265'
$ws.Range("E16").Value = '  -2.25%  '
$ws.Range("D17").Value = '1.628.16'
$ws.Range("E17").Value = '  +1.45%  '
$ws.Range("D18").Value = '''94.25'
$ws.Range("E18").Value = '  +0.24%  '
$ws.Range("D19").Value = '''0.06868'
$ws.Range("E19").Value = '  -0.48%  '
$ws.Range("D20").Value = '''18.07'
$ws.Range("E20").Value = '  -0.63%  '
$ws.Range("D21").Value = '''6.504'
$ws.Range("E21").Value = '  -0.53%  '
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("D23").Value = '23.403.59'
$ws.Range("E23").Value = '  +0.72%  '
$ws.Range("D24").Value = '''12.72'
$ws.Range("E24").Value = '  -1.90%  '
$ws.Range("D25").Value = '''2.416'
$ws.Range("E25").Value = '  +1.30%  '
$ws.Range("D26").Value = '''3.011'
$ws.Range("E26").Value = '  -0.25%  '
$ws.Range("D27").Value = '''21.11'
$ws.Range("E27").Value = '  -0.66%  '
$ws.Range("D28").Value = '''151.01'
$ws.Range("E28").Value = '  +0.34%  '
$ws.Range("D29").Value = '''5.327'
$ws.Range("E29").Value = '  +1.30%  '
$ws.Range("D30").Value = '''134.58'
$ws.Range("E30").Value = '  +0.29%  '
$ws.Range("D31").Value = '''2.289'
$ws.Range("E31").Value = '  -4.01%  '
$ws.Range("D32").Value = '1.807.86'
$ws.Range("E32").Value = '  +1.37%  '
$ws.Range("D33").Value = '''6.739'
$ws.Range("E33").Value = '  -0.57%  '
$ws.Range("D34").Value = '''0.9458'
$ws.Range("E34").Value = '  -2.12%  '
$ws.Range("D35").Value = '''0.02820'
$ws.Range("E35").Value = '  +2.61%  '
$ws.Range("D36").Value = '''10.29'
$ws.Range("E36").Value = '  +0.22%  '
$ws.Range("D37").Value = '''0.2517'
$ws.Range("E37").Value = '  -0.82%  '
$ws.Range("D38").Value = '''0.07185'
$ws.Range("E38").Value = '  -4.74%  '
$ws.Range("D39").Value = '''0.08752'
$ws.Range("E39").Value = '  -0.73%  '
$ws.Range("D40").Value = '''6.044'
$ws.Range("E40").Value = '  -1.30%  '
$ws.Range("D41").Value = '''1.367'
$ws.Range("E41").Value = '  -2.03%  '
$ws.Range("D42").Value = '''0.7015'
$ws.Range("E42").Value = '  -1.60%  '
$ws.Range("D43").Value = '''12.39'
$ws.Range("E43").Value = '  -1.43%  '
$ws.Range("D44").Value = '''15.91'
$ws.Range("E44").Value = '  +1.65%  '
$ws.Range("D45").Value = '''0.6480'
$ws.Range("E45").Value = '  -1.12%  '
$ws.Range("D46").Value = '''2.318'
$ws.Range("E46").Value = '  -0.56%  '
$ws.Range("D47").Value = '''0.9998'
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("D48").Value = '''4.003'
$ws.Range("E48").Value = '  -0.33%  '
$ws.Range("D49").Value = '''0.07959'
$ws.Range("E49").Value = '  +0.06%  '
$ws.Range("D50").Value = '''128.23'
$ws.Range("E50").Value = '  -3.38%  '
$ws.Range("D51").Value = '''1.191'
$ws.Range("E51").Value = '  -1.19%  '
